$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.975.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.750.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -16.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.236.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.661.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.756.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.542"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0901"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.980"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.14%  "
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "325.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0590"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0255"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
